$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.804
$ws.Range("L2").Value = 0.499
$ws.Range("G3").Value = 0.8080000000000001
$ws.Range("G4").Value = 0.734
$ws.Range("L4").Value = 0.457
$ws.Range("L5").Value = 0.429
$ws.Range("C6").Value = 0.661
$ws.Range("E6").Value = 0.598
$ws.Range("E7").Value = 0.577
$ws.Range("K7").Value = 0.481
$ws.Range("D8").Value = 0.487
$ws.Range("E8").Value = 0.596
$ws.Range("C9").Value = 0.663
$ws.Range("E9").Value = 0.608
$ws.Range("G10").Value = 0.73
$ws.Range("I10").Value = 0.544
$ws.Range("C11").Value = 0.6909999999999999
$ws.Range("L11").Value = 0.443
$ws.Range("J13").Value = 0.384
$ws.Range("E14").Value = 0.578
$ws.Range("L16").Value = 0.43
$ws.Range("G17").Value = 0.742
$ws.Range("E21").Value = 0.633
$ws.Range("C23").Value = 0.671
$ws.Range("L23").Value = 0.446
$ws.Range("C32").Value = 0.571
$ws.Range("H43").Value = 0.431
$ws.Range("C46").Value = 0.628
$ws.Range("G49").Value = 0.599
$ws.Range("H56").Value = 0.434
$ws.Range("I56").Value = 0.537
$ws.Range("L58").Value = 0.431
$ws.Range("H62").Value = 0.424
$ws.Range("I62").Value = 0.591
$ws.Range("J65").Value = 0.379
$ws.Range("G68").Value = 0.632
$ws.Range("I68").Value = 0.546
$ws.Range("D69").Value = 0.538
$ws.Range("B71").Value = 0.649
$ws.Range("E71").Value = 0.584
$ws.Range("F77").Value = 0.527
$ws.Range("H77").Value = 0.376
$ws.Range("J79").Value = 0.4
$ws.Range("F80").Value = 0.537
$ws.Range("G81").Value = 0.618
$ws.Range("B82").Value = 0.667
$ws.Range("C82").Value = 0.675
$ws.Range("J84").Value = 0.404
$ws.Range("E87").Value = 0.5620000000000001
$ws.Range("H87").Value = 0.374
$ws.Range("I88").Value = 0.638
$ws.Range("G89").Value = 0.671
$ws.Range("J90").Value = 0.366
$ws.Range("B91").Value = 0.632
$ws.Range("H91").Value = 0.388
$ws.Range("H93").Value = 0.464
$ws.Range("C95").Value = 0.5610000000000001
$ws.Range("L95").Value = 0.386
$ws.Range("G97").Value = 0.675
$ws.Range("L97").Value = 0.46
$ws.Range("B98").Value = 0.587
$ws.Range("K98").Value = 0.414
$ws.Range("C99").Value = 0.649
$ws.Range("H99").Value = 0.392
$ws.Range("C100").Value = 0.571
$ws.Range("E100").Value = 0.515
